# Applies the "Incorporate addenda a, b, and c" change:
#  1. Bump schema_version on the RS0004 sheet from 1.0.0 to 2.0.0
#  2. Relax the three list data validations (no input/error messages)
#  3. Add a new grid/lookup column J ("operation_state") to the
#     performance_map_cooling sheet, with header, units and per-row
#     "NORMAL" values, plus an explanatory cell comment on J3.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) RS0004 sheet: schema_version bump + data validation tweaks
# ---------------------------------------------------------------------
$rs = $wb.Worksheets.Item("RS0004")

$rs.Range("C4").Value2 = "2.0.0"

foreach ($addr in @("C6", "C20", "C22")) {
    $v = $rs.Range($addr).Validation
    $v.ShowInput = $false
    $v.ShowError = $false
}

# ---------------------------------------------------------------------
# 2) performance_map_cooling sheet: new column J (operation_state)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("performance_map_cooling")

# xlPasteFormats - used instead of Range.Style so the cell-level direct
# formatting (e.g. the 45-degree header rotation) carries over, not just
# the underlying named cell style.
$xlPasteFormats = -4122

# Row 2: part of the "lookup_variables" group header/banner -> blank
# cell styled like its neighbours (H2 / I2).
$ws.Range("I2").Copy()
$ws.Range("J2").PasteSpecial($xlPasteFormats)

# Row 3: column name, styled like the other lookup variable headers.
$ws.Range("I3").Copy()
$ws.Range("J3").PasteSpecial($xlPasteFormats)
$ws.Range("J3").Value2 = "operation_state"

# Row 4: units ("-"), styled like the other unit cells.
$ws.Range("I4").Copy()
$ws.Range("J4").PasteSpecial($xlPasteFormats)
$ws.Range("J4").Value2 = "-"

# Rows 5-68: data values, all "NORMAL", styled like the rest of the row.
for ($r = 5; $r -le 68; $r++) {
    $ws.Cells.Item($r, 9).Copy()
    $ws.Cells.Item($r, 10).PasteSpecial($xlPasteFormats)
    $ws.Cells.Item($r, 10).Value2 = "NORMAL"
}

$excel.CutCopyMode = 0

# New comment explaining the added column.
$ws.Range("J3").AddComment("The operation state at the operating conditions")
